# Auto-generated Excel COM-interop script applying the Titan_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

# Sheet ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 770.2857
$ws.Range("I18").Value = 530
$ws.Range("K18").Value = 530
$ws.Range("M18").Value = -246

# Sheet ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 11111780
$ws.Range("I100").Value = 16667272
$ws.Range("J100").Value = 795.2
$ws.Range("K100").Value = 16667272
$ws.Range("L100").Value = 795.2
$ws.Range("M100").Value = -16666731
$ws.Range("N100").Value = -1877.2

# Sheet ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5324127.5
$ws.Range("I116").Value = 5536972.5
$ws.Range("K116").Value = 5536972.5
$ws.Range("M116").Value = -5533530.5

# Sheet ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 265057.34
$ws.Range("I132").Value = 312326.88
$ws.Range("J132").Value = 60222.668
$ws.Range("K132").Value = 936980.64
$ws.Range("L132").Value = 180668.004
$ws.Range("M132").Value = -934450.64
$ws.Range("N132").Value = -185728.004

# Sheet ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 16800
$ws.Range("J133").Value = 16800
$ws.Range("L133").Value = 16800
$ws.Range("N133").Value = -26920

# Sheet ALC row 134
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 55000
$ws.Range("J134").Value = 55000
$ws.Range("L134").Value = 55000
$ws.Range("N134").Value = -65140

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 23810666
$ws.Range("I137").Value = 26316652
$ws.Range("J137").Value = 3800.75
$ws.Range("K137").Value = 78949956
$ws.Range("L137").Value = 11402.25
$ws.Range("M137").Value = -78947406
$ws.Range("N137").Value = -16502.25

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16965.94
$ws.Range("I32").Value = 2111.1746
$ws.Range("K32").Value = 2111.1746
$ws.Range("M32").Value = -1824.1746

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1216.25
$ws.Range("I45").Value = 926.6667
$ws.Range("K45").Value = 926.6667
$ws.Range("M45").Value = -549.6667

# Sheet ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2292.625
$ws.Range("I61").Value = 1568.1765
$ws.Range("J61").Value = 4052
$ws.Range("K61").Value = 1568.1765
$ws.Range("L61").Value = 4052
$ws.Range("M61").Value = -1356.1765
$ws.Range("N61").Value = -4476

# Sheet ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 59261
$ws.Range("J133").Value = 59261
$ws.Range("L133").Value = 59261
$ws.Range("N133").Value = -64321

# Sheet ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2292.625
$ws.Range("I136").Value = 1568.1765
$ws.Range("J136").Value = 4052
$ws.Range("K136").Value = 4704.529500000001
$ws.Range("L136").Value = 12156
$ws.Range("M136").Value = -2154.529500000001
$ws.Range("N136").Value = -17256

# Sheet ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 44085.832
$ws.Range("J139").Value = 44085.832
$ws.Range("L139").Value = 44085.832
$ws.Range("N139").Value = -54365.832

# Sheet BSM row 59
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 42816.668
$ws.Range("J59").Value = 47380
$ws.Range("L59").Value = 47380
$ws.Range("N59").Value = -49074

# Sheet BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4532.6665
$ws.Range("I86").Value = 1700.9412
$ws.Range("J86").Value = 7541.375
$ws.Range("K86").Value = 1700.9412
$ws.Range("L86").Value = 7541.375
$ws.Range("M86").Value = -577.9412
$ws.Range("N86").Value = -9787.375

# Sheet BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4532.6665
$ws.Range("I89").Value = 1700.9412
$ws.Range("J89").Value = 7541.375
$ws.Range("K89").Value = 8504.706
$ws.Range("L89").Value = 37706.875
$ws.Range("M89").Value = -2888.706
$ws.Range("N89").Value = -48938.875

# Sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4115.4287
$ws.Range("I134").Value = 2161.2
$ws.Range("J134").Value = 5892
$ws.Range("K134").Value = 6483.599999999999
$ws.Range("L134").Value = 17676
$ws.Range("M134").Value = -3948.599999999999
$ws.Range("N134").Value = -22746

# Sheet CRP row 21
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 9500
$ws.Range("I21").Value = 9000
$ws.Range("K21").Value = 9000
$ws.Range("M21").Value = -8765

# Sheet CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 539.26666
$ws.Range("I22").Value = 380.66666
$ws.Range("J22").Value = 777.1667
$ws.Range("K22").Value = 380.66666
$ws.Range("L22").Value = 777.1667
$ws.Range("M22").Value = -30.66665999999998
$ws.Range("N22").Value = -1477.1667

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1345.6888
$ws.Range("I31").Value = 883.7353000000001
$ws.Range("J31").Value = 2773.5454
$ws.Range("K31").Value = 883.7353000000001
$ws.Range("L31").Value = 2773.5454
$ws.Range("M31").Value = -588.7353000000001
$ws.Range("N31").Value = -3363.5454

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1345.6888
$ws.Range("I34").Value = 883.7353000000001
$ws.Range("J34").Value = 2773.5454
$ws.Range("K34").Value = 883.7353000000001
$ws.Range("L34").Value = 2773.5454
$ws.Range("M34").Value = -681.7353000000001
$ws.Range("N34").Value = -3177.5454

# Sheet CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1759.6875
$ws.Range("I94").Value = 1101.8334
$ws.Range("J94").Value = 2154.4
$ws.Range("K94").Value = 1101.8334
$ws.Range("L94").Value = 2154.4
$ws.Range("M94").Value = -650.8334
$ws.Range("N94").Value = -3056.4

# Sheet GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 53.5
$ws.Range("I2").Value = 51
$ws.Range("K2").Value = 51
$ws.Range("M2").Value = 62

# Sheet GSM row 17
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 6215.9
$ws.Range("I17").Value = 10400
$ws.Range("J17").Value = 2031.8
$ws.Range("K17").Value = 10400
$ws.Range("L17").Value = 2031.8
$ws.Range("M17").Value = -10232
$ws.Range("N17").Value = -2367.8

# Sheet GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6618.909
$ws.Range("I70").Value = 6901
$ws.Range("K70").Value = 6901
$ws.Range("M70").Value = -6631

# Sheet GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6618.909
$ws.Range("I73").Value = 6901
$ws.Range("K73").Value = 6901
$ws.Range("M73").Value = -5965

# Sheet GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2190
$ws.Range("I80").Value = 2100
$ws.Range("J80").Value = 2212.5
$ws.Range("K80").Value = 2100
$ws.Range("L80").Value = 2212.5
$ws.Range("M80").Value = -1102
$ws.Range("N80").Value = -4208.5

# Sheet GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2190
$ws.Range("I83").Value = 2100
$ws.Range("J83").Value = 2212.5
$ws.Range("K83").Value = 10500
$ws.Range("L83").Value = 11062.5
$ws.Range("M83").Value = -5508
$ws.Range("N83").Value = -21046.5

# Sheet GSM row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 64425
$ws.Range("J138").Value = 64425
$ws.Range("L138").Value = 64425
$ws.Range("N138").Value = -74705

# Sheet LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 597.36365
$ws.Range("I16").Value = 597.36365
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 597.36365
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = -427.36365
$ws.Range("M16").Value = ""

# Sheet LTW row 20
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 4500
$ws.Range("J20").Value = 4500
$ws.Range("L20").Value = 4500
$ws.Range("N20").Value = -4952

# Sheet LTW row 24
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 2564.25
$ws.Range("J24").Value = 2564.25
$ws.Range("L24").Value = 2564.25
$ws.Range("N24").Value = -3250.25

# Sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5147.6665
$ws.Range("I132").Value = 4610
$ws.Range("K132").Value = 13830
$ws.Range("M132").Value = -11300

# Sheet WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 44492.957
$ws.Range("I122").Value = 73557.78999999999
$ws.Range("J122").Value = 3802.2
$ws.Range("K122").Value = 220673.37
$ws.Range("L122").Value = 11406.6
$ws.Range("M122").Value = -218223.37
$ws.Range("N122").Value = -16306.6

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 23814970
$ws.Range("I132").Value = 62510496
$ws.Range("J132").Value = 2336.1538
$ws.Range("K132").Value = 187531488
$ws.Range("L132").Value = 7008.4614
$ws.Range("M132").Value = -187528958
$ws.Range("N132").Value = -12068.4614

# Sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9835109
$ws.Range("I136").Value = 12860125
$ws.Range("J136").Value = 3806.375
$ws.Range("K136").Value = 38580375
$ws.Range("L136").Value = 11419.125
$ws.Range("M136").Value = -38577825
$ws.Range("N136").Value = -16519.125

Write-Host "Applied all cell updates"
